# A new weekly price record needs to be inserted as row 137 (pushing the
# existing rows 137-228 down to 138-229). The new record carries the same
# constant columns (Mercado ID, Mercado, Region, Codreg, Tipo, Producto ID,
# Producto, Categoria ID, Categoria, Variedad) as the rest of the block, and
# fresh values for the remaining (varying) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 137; everything currently at/after
# row 137 (up to the last used row, 228) shifts down to 138..229.
$ws.Rows.Item(137).Insert()

# Columns that stay constant across this whole data block - copy them down
# from the row immediately above (row 136) into the freshly inserted row 137.
for ($col = 1; $col -le 11; $col++) {
    if ($col -ne 4) {
        $ws.Cells.Item(137, $col).Value = $ws.Cells.Item(136, $col).Value2
    }
}

# New record's own data (columns D, L-T).
$ws.Cells.Item(137, 4).Value  = 44827                                      # D137 Fecha
$ws.Cells.Item(137, 12).Value = "Primera"                                  # L137 Calidad
$ws.Cells.Item(137, 13).Value = 40                                         # M137 Volumen
$ws.Cells.Item(137, 14).Value = 12000                                      # N137 Precio minimo
$ws.Cells.Item(137, 15).Value = 12000                                      # O137 Precio maximo
$ws.Cells.Item(137, 16).Value = 12000                                      # P137 Precio promedio ponderado
$ws.Cells.Item(137, 17).Value = "$/bandeja 12 canastillos 125 gramos"      # Q137 Unidad de comercializacion
$ws.Cells.Item(137, 18).Value = "Provincia de Limarí"                      # R137 Origen
$ws.Cells.Item(137, 19).Value = 8000                                       # S137 Precio $/Kg
$ws.Cells.Item(137, 20).Value = 1.5                                        # T137 Kg / unidad

# Match the date-time number format used by the rest of column D.
$ws.Cells.Item(137, 4).NumberFormat = $ws.Cells.Item(136, 4).NumberFormat
